$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: delete all comments. This removes commentRangeStart/End and
# commentReference runs from the body, and empties the comments parts.
# ---------------------------------------------------------------------
while ($d.Comments.Count -gt 0) {
    $d.Comments.Item(1).Delete()
}

# ---------------------------------------------------------------------
# Step 2: paragraph with the Lexico definition quote - replace the
# trailing "(Add Reference)" run with a properly formatted citation,
# including the spell-check markers around "Lexico".
# ---------------------------------------------------------------------
$p7 = $d.Paragraphs.Item(7)
$rng7 = $p7.Range.Duplicate
$rng7.Find.Execute("(Add Reference)", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng7.Text = "X"
$lexicoXml = @'
<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="21"/></w:rPr><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="21"/></w:rPr><w:t>Lexico</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="21"/></w:rPr><w:t xml:space="preserve"> Dictionaries | English, 2019</w:t></w:r><w:r><w:rPr><w:sz w:val="21"/></w:rPr><w:t>)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rng7.InsertXML($lexicoXml)
$rng7.Text = ""

# ---------------------------------------------------------------------
# Step 3: remove the trailing space after "your data/systems."
# ---------------------------------------------------------------------
$d.Content.Find.Execute("your data/systems. ", $false, $false, $false, $false, $false, $true, 1, $false, "your data/systems.", 2)

# ---------------------------------------------------------------------
# Step 4: PSN attack paragraph - swap the three "(Add reference)" /
# "(Add Reference)" placeholders for the real citations.
# ---------------------------------------------------------------------
$p14 = $d.Paragraphs.Item(14)

$rngA = $p14.Range.Duplicate
$rngA.Find.Execute("(Add reference)", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngA.Text = "(Sony Global - Sony Global Headquarters, 2019)"

$rngB = $p14.Range.Duplicate
$rngB.Find.Execute("(Add reference)", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngB.Text = "(Cbsnews.com, 2019)"

$rngC = $p14.Range.Duplicate
$rngC.Find.Execute("(Add Reference). ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngC.Text = "(Seek Market Insights AU, 2019)"

# ---------------------------------------------------------------------
# Step 5: insert a new, empty ListParagraph-styled paragraph right
# after the PSN paragraph (before the pre-existing blank paragraph).
# ---------------------------------------------------------------------
$endRng14 = $p14.Range.Duplicate
$endRng14.Collapse(0)
$endRng14.InsertParagraphAfter()
$p15 = $d.Paragraphs.Item(15)
$cleanRng15 = $p15.Range.Duplicate
$cleanRng15.Collapse(1)
$emptyListParaXml = @'
<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="1080"/></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$cleanRng15.InsertXML($emptyListParaXml)

# ---------------------------------------------------------------------
# Step 6: append a new paragraph "(Cbsnews.com, 2019)" right after the
# "...tailored assistance to the elderly." paragraph, before the final
# empty paragraph.
# ---------------------------------------------------------------------
$lastTextParaIndex = $d.Paragraphs.Count - 1
$pElderly = $d.Paragraphs.Item($lastTextParaIndex)
$endRngElderly = $pElderly.Range.Duplicate
$endRngElderly.Collapse(0)
$endRngElderly.InsertParagraphAfter()
$pFinalEmpty = $d.Paragraphs.Item($lastTextParaIndex + 1)
$cleanRngFinal = $pFinalEmpty.Range.Duplicate
$cleanRngFinal.Collapse(1)
$cbsXml = @'
<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>(Cbsnews.com, 2019)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$cleanRngFinal.InsertXML($cbsXml)

Write-Host "Done."
